$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 17 corresponds to "Chapter 15: Arrays" - mark it as DONE like the rows above it.
$ws.Range("B17").Value = "DONE"
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial(-4122) # xlPasteFormats
